$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1) by copying its format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for column I (I0)
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 3
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 5

# New data values for column J (IF)
$ws.Range("J2").Value = 5
$ws.Range("J3").Value = 9
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 9
